$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price(D) and Volume(E) columns to remain text, not get auto-converted to numbers
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.045.19"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.853.51"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").Value = "1.014"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "310.48"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "0.4781"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("D8").Value = "0.3682"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").Value = "0.07233"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "0.9315"
$ws.Range("E10").Value = "  +3.23%  "
$ws.Range("D11").Value = "19.74"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "0.07737"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.871.67"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").Value = "5.332"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "6.437"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "89.03"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "1.016"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "0.000008647"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "1.013"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "27.057.91"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "14.56"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").Value = "5.063"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "10.68"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "1.933"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "152.94"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "18.23"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("D27").Value = "2.008"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").Value = "114.50"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").Value = "4.998"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").Value = "0.08905"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").Value = "3.332"
$ws.Range("E31").Value = "  +5.79%  "
$ws.Range("D32").Value = "1.178"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "0.7459"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "4.509"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").Value = "2.733"
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("D36").Value = "1.112"
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("D37").Value = "0.05279"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").Value = "0.01956"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").Value = "2.978"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").Value = "0.5217"
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("D41").Value = "7.028"
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").Value = "0.1512"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").Value = "8.228"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("D45").Value = "0.4747"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("D46").Value = "1.016"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "101.76"
$ws.Range("E47").Value = "  +3.86%  "
$ws.Range("D48").Value = "1.611"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").Value = "65.78"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").Value = "0.06032"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").Value = "0.8892"
$ws.Range("E51").Value = "  +4.30%  "

# Restore default styling (remove the temporary text number format) without altering values
$ws.Range("D2:E51").Style = "Normal"

